$d = $word.ActiveDocument

$pairs = @(
    @("20×87=", "24×85="),
    @("27×48=", "13×83="),
    @("75×41=", "95×33="),
    @("31×78=", "89×81="),
    @("91×49=", "83×93="),
    @("95×97=", "63×37="),
    @("31×85=", "56×38="),
    @("42×45=", "95×94="),
    @("76×31=", "25×97="),
    @("31×36=", "86×83="),
    @("37×32=", "98×95="),
    @("27×69=", "26×97="),
    @("76×78=", "63×60="),
    @("34×90=", "11×24="),
    @("12×47=", "47×21="),
    @("82×81=", "65×26="),
    @("72×20=", "41×20="),
    @("76×24=", "93×24="),
    @("71×53=", "63×26="),
    @("80×75=", "79×60="),
    @("82×91=", "32×13="),
    @("41×14=", "39×14="),
    @("24×23=", "64×79="),
    @("19×67=", "56×19="),
    @("79×76=", "54×89=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
